$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.43775051914092
$ws.Range("C2").Value = 11.19696365860557
$ws.Range("E2").Value = 15.81573047652792
$ws.Range("F2").Value = 38.38584959282824
$ws.Range("G2").Value = 3.653661470515581
$ws.Range("J2").Value = 8.045255017942342
$ws.Range("L2").Value = 12.34078765958234
$ws.Range("M2").Value = 16.99452758298961
$ws.Range("N2").Value = 18.92523607500412
$ws.Range("O2").Value = 24.89179599123506
$ws.Range("B3").Value = 16.05685461848186
$ws.Range("C3").Value = 11.08167308098274
$ws.Range("E3").Value = 15.84848484082278
$ws.Range("F3").Value = 38.42198615855883
$ws.Range("G3").Value = 3.655765374260229
$ws.Range("J3").Value = 8.04355074766204
$ws.Range("L3").Value = 12.3393455301605
$ws.Range("M3").Value = 16.91562405843424
$ws.Range("N3").Value = 18.98375408561752
$ws.Range("O3").Value = 24.95110627777614
$ws.Range("B4").Value = 15.82103685662449
$ws.Range("C4").Value = 11.00956860142087
$ws.Range("E4").Value = 15.87010679756442
$ws.Range("F4").Value = 38.45306125204628
$ws.Range("G4").Value = 3.657126323112413
$ws.Range("J4").Value = 8.042618807117867
$ws.Range("L4").Value = 12.33985639173106
$ws.Range("M4").Value = 16.86937108954217
$ws.Range("N4").Value = 19.02152684281691
$ws.Range("O4").Value = 24.99330753928308
$ws.Range("B5").Value = 15.7245895945912
$ws.Range("C5").Value = 10.97987005328089
$ws.Range("E5").Value = 15.8792984263064
$ws.Range("F5").Value = 38.46795725852268
$ws.Range("G5").Value = 3.657698362647787
$ws.Range("J5").Value = 8.042268027620802
$ws.Range("L5").Value = 12.34041699036301
$ws.Range("M5").Value = 16.85108789673904
$ws.Range("N5").Value = 19.03738405197098
$ws.Range("O5").Value = 25.01195582438472
$ws.Range("B6").Value = 15.70855767230491
$ws.Range("C6").Value = 10.97492003812748
$ws.Range("E6").Value = 15.88084768974583
$ws.Range("F6").Value = 38.47056549694156
$ws.Range("G6").Value = 3.657794404443713
$ws.Range("J6").Value = 8.042211538078945
$ws.Range("L6").Value = 12.34053139866096
$ws.Range("M6").Value = 16.84808650992975
$ws.Range("N6").Value = 19.04004521988999
$ws.Range("O6").Value = 25.01513987837771
$ws.Range("B7").Value = 15.81973736212353
$ws.Range("C7").Value = 11.00916933287228
$ws.Range("E7").Value = 15.87022921757041
$ws.Range("F7").Value = 38.45325310832485
$ws.Range("G7").Value = 3.657133967143163
$ws.Range("J7").Value = 8.042613958722017
$ws.Range("L7").Value = 12.33986252358235
$ws.Range("M7").Value = 16.8691222098766
$ws.Range("N7").Value = 19.02173881633207
$ws.Range("O7").Value = 24.99355316687283
$ws.Range("B8").Value = 16.30690367560334
$ws.Range("C8").Value = 11.15749370466731
$ws.Range("E8").Value = 15.82671103342965
$ws.Range("F8").Value = 38.39646398106763
$ws.Range("G8").Value = 3.654372579050091
$ws.Range("J8").Value = 8.044643745904985
$ws.Range("L8").Value = 12.3400014081068
$ws.Range("M8").Value = 16.96687457461605
$ws.Range("N8").Value = 18.94503140574794
$ws.Range("O8").Value = 24.91104350664853
$ws.Range("B9").Value = 17.24094939588184
$ws.Range("C9").Value = 11.4372557040078
$ws.Range("E9").Value = 15.75332962579535
$ws.Range("F9").Value = 38.35567261449398
$ws.Range("G9").Value = 3.649503612430664
$ws.Range("J9").Value = 8.049523880070227
$ws.Range("L9").Value = 12.35129349869057
$ws.Range("M9").Value = 17.17535377823447
$ws.Range("N9").Value = 18.80917245890345
$ws.Range("O9").Value = 24.79529073802977
$ws.Range("B10").Value = 17.90686331029642
$ws.Range("C10").Value = 11.63516469556863
$ws.Range("E10").Value = 15.70666675035853
$ws.Range("F10").Value = 38.36874986560609
$ws.Range("G10").Value = 3.646255785185396
$ws.Range("J10").Value = 8.05364659930077
$ws.Range("L10").Value = 12.36622408989316
$ws.Range("M10").Value = 17.3378870565553
$ws.Range("N10").Value = 18.71815816259714
$ws.Range("O10").Value = 24.73850857102893
$ws.Range("B11").Value = 18.20398318627533
$ws.Range("C11").Value = 11.72336449800891
$ws.Range("E11").Value = 15.68700466291326
$ws.Range("F11").Value = 38.3840306296927
$ws.Range("G11").Value = 3.644849042785489
$ws.Range("J11").Value = 8.055636408728466
$ws.Range("L11").Value = 12.37443722320342
$ws.Range("M11").Value = 17.4136647489029
$ws.Range("N11").Value = 18.67864822765579
$ws.Range("O11").Value = 24.71884702555958
$ws.Range("B12").Value = 18.31555304911257
$ws.Range("C12").Value = 11.75648585938785
$ws.Range("E12").Value = 15.67978355746653
$ws.Range("F12").Value = 38.39115575142755
$ws.Range("G12").Value = 3.644326457375307
$ws.Range("J12").Value = 8.056406128544596
$ws.Range("L12").Value = 12.3777498078387
$ws.Range("M12").Value = 17.44260718592488
$ws.Range("N12").Value = 18.66395785489755
$ws.Range("O12").Value = 24.71229088751926
$ws.Range("B13").Value = 18.29156805320313
$ws.Range("C13").Value = 11.74936517695057
$ws.Range("E13").Value = 15.68132877587284
$ws.Range("F13").Value = 38.38956175700444
$ws.Range("G13").Value = 3.644438556321637
$ws.Range("J13").Value = 8.056239637698557
$ws.Range("L13").Value = 12.37702741095437
$ws.Range("M13").Value = 17.43636320966365
$ws.Range("N13").Value = 18.6671096446213
$ws.Range("O13").Value = 24.71366328861701
$ws.Range("B14").Value = 18.21318159371285
$ws.Range("C14").Value = 11.72609506713582
$ws.Range("E14").Value = 15.68640608208676
$ws.Range("F14").Value = 38.38459001232393
$ws.Range("G14").Value = 3.644805846890272
$ws.Range("J14").Value = 8.055699410210245
$ws.Range("L14").Value = 12.3747057054426
$ws.Range("M14").Value = 17.41604101100664
$ws.Range("N14").Value = 18.67743421271696
$ws.Range("O14").Value = 24.71828981496688
$ws.Range("B15").Value = 18.16504175736792
$ws.Range("C15").Value = 11.71180481743044
$ws.Range("E15").Value = 15.68954529880756
$ws.Range("F15").Value = 38.38171888800836
$ws.Range("G15").Value = 3.645032139057011
$ws.Range("J15").Value = 8.055370611222303
$ws.Range("L15").Value = 12.37330990027633
$ws.Range("M15").Value = 17.40362472504333
$ws.Range("N15").Value = 18.68379359482313
$ws.Range("O15").Value = 24.72123956185533
$ws.Range("B16").Value = 17.88731955733273
$ws.Range("C16").Value = 11.62936252876882
$ws.Range("E16").Value = 15.70798315989806
$ws.Range("F16").Value = 38.36793875693838
$ws.Range("G16").Value = 3.64634913764413
$ws.Range("J16").Value = 8.053518842830702
$ws.Range("L16").Value = 12.36571577840416
$ws.Range("M16").Value = 17.33297034628227
$ws.Range("N16").Value = 18.72077823277077
$ws.Range("O16").Value = 24.7399178293778
$ws.Range("B17").Value = 17.71537917980374
$ws.Range("C17").Value = 11.57830745868844
$ws.Range("E17").Value = 15.71969464991286
$ws.Range("F17").Value = 38.36187321651312
$ws.Range("G17").Value = 3.647175148089675
$ws.Range("J17").Value = 8.052411976019977
$ws.Range("L17").Value = 12.36141969209386
$ws.Range("M17").Value = 17.29008517146414
$ws.Range("N17").Value = 18.74395120308097
$ws.Range("O17").Value = 24.7529580316816
$ws.Range("B18").Value = 17.61594290705834
$ws.Range("C18").Value = 11.54877054700974
$ws.Range("E18").Value = 15.72657813598448
$ws.Range("F18").Value = 38.35926310195404
$ws.Range("G18").Value = 3.647656906061414
$ws.Range("J18").Value = 8.051786113421212
$ws.Range("L18").Value = 12.35908258129619
$ws.Range("M18").Value = 17.26559306688766
$ws.Range("N18").Value = 18.75745791335547
$ws.Range("O18").Value = 24.76103895650179
$ws.Range("B19").Value = 17.58218604450215
$ws.Range("C19").Value = 11.53874087963166
$ws.Range("E19").Value = 15.72893409005903
$ws.Range("F19").Value = 38.35853034866994
$ws.Range("G19").Value = 3.647821166207098
$ws.Range("J19").Value = 8.051576064989769
$ws.Range("L19").Value = 12.35831432695419
$ws.Range("M19").Value = 17.25733092949123
$ws.Range("N19").Value = 18.76206169602189
$ws.Range("O19").Value = 24.76387465521159
$ws.Range("B20").Value = 17.73373929549464
$ws.Range("C20").Value = 11.58376020468456
$ws.Range("E20").Value = 15.71843269602012
$ws.Range("F20").Value = 38.36242798966246
$ws.Range("G20").Value = 3.647086529094258
$ws.Range("J20").Value = 8.052528690014782
$ws.Range("L20").Value = 12.36186317477618
$ws.Range("M20").Value = 17.29463245696236
$ws.Range("N20").Value = 18.74146596205449
$ws.Range("O20").Value = 24.75150977856079
$ws.Range("B21").Value = 18.23623200851768
$ws.Range("C21").Value = 11.73293772314843
$ws.Range("E21").Value = 15.68490866573519
$ws.Range("F21").Value = 38.38601403583836
$ws.Range("G21").Value = 3.644697690564981
$ws.Range("J21").Value = 8.055857649642112
$ws.Range("L21").Value = 12.37538216792931
$ws.Range("M21").Value = 17.42200356994048
$ws.Range("N21").Value = 18.67439428500832
$ws.Range("O21").Value = 24.71690674236307
$ws.Range("B22").Value = 18.55909930010667
$ws.Range("C22").Value = 11.82880694492969
$ws.Range("E22").Value = 15.66430707646093
$ws.Range("F22").Value = 38.40922913509023
$ws.Range("G22").Value = 3.643195395245333
$ws.Range("J22").Value = 8.058127767484951
$ws.Range("L22").Value = 12.38539668348231
$ws.Range("M22").Value = 17.5066796661488
$ws.Range("N22").Value = 18.63213926511417
$ws.Range("O22").Value = 24.69947568209303
$ws.Range("B23").Value = 18.38731945986049
$ws.Range("C23").Value = 11.77779338924345
$ws.Range("E23").Value = 15.6751830063021
$ws.Range("F23").Value = 38.39612644234022
$ws.Range("G23").Value = 3.643991821355734
$ws.Range("J23").Value = 8.056907594913296
$ws.Range("L23").Value = 12.37994451969817
$ws.Range("M23").Value = 17.46136133938484
$ws.Range("N23").Value = 18.65454730758062
$ws.Range("O23").Value = 24.70830401464396
$ws.Range("B24").Value = 17.72544050245827
$ws.Range("C24").Value = 11.58129559233375
$ws.Range("E24").Value = 15.71900275687837
$ws.Range("F24").Value = 38.36217444453769
$ws.Range("G24").Value = 3.647126572327708
$ws.Range("J24").Value = 8.052475890918922
$ws.Range("L24").Value = 12.36166226256296
$ws.Range("M24").Value = 17.29257611963686
$ws.Range("N24").Value = 18.7425889654664
$ws.Range("O24").Value = 24.75216271500284
$ws.Range("B25").Value = 16.99133158185014
$ws.Range("C25").Value = 11.36285158976663
$ws.Range("E25").Value = 15.77190516067991
$ws.Range("F25").Value = 38.35914601428431
$ws.Range("G25").Value = 3.650762697642334
$ws.Range("J25").Value = 8.048108503555207
$ws.Range("L25").Value = 12.34706718153707
$ws.Range("M25").Value = 17.11724509390014
$ws.Range("N25").Value = 18.84437479065066
$ws.Range("O25").Value = 24.82165423272212
